{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Updates the stack-trace text in the document (a Java exception dump\n// embedded as plain text in a single run) to reflect refactored line\n// numbers / stack frames, per a M2DocEvaluator 3.1.1 -> 3.2.0 code move:\n//  - several M2DocEvaluator.java / M2DocUtils.java / AbstractTemplatesTestSuite.java\n//    line numbers shift,\n//  - a generated accessor class name changes,\n//  - the tail of the trace (Eclipse JDT JUnit launcher frames) is replaced\n//    by a Maven Surefire / Tycho / Equinox launcher frame sequence.\nconst REPLACEMENTS = [\n  [\"org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1207)\", \"org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1296)\", 1],\n  [\"org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)\", \"org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)\", 3],\n  [\"org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1467)\", \"org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1556)\", 1],\n  [\"org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:297)\", \"org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:301)\", 1],\n  [\"org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:282)\", \"org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:286)\", 1],\n  [\"org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:845)\", \"org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)\", 1],\n  [\"org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514)\", \"org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)\", 1],\n  [\"org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421)\", \"org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)\", 1],\n  [\"sun.reflect.GeneratedMethodAccessor73.invoke(Unknown Source)\", \"sun.reflect.GeneratedMethodAccessor5.invoke(Unknown Source)\", 1],\n  [\"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\", \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)\\n\\tat org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:161)\\n\\tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)\\n\\tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)\\n\\tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)\\n\\tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)\\n\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)\\n\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)\\n\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)\\n\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)\\n\\tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)\\n\\tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)\\n\\tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)\", 1],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText, expectedCount] of REPLACEMENTS) {\n  const results = body.search(oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== expectedCount) {\n    throw new Error(\n      \"Unexpected match count for replacement (expected \" + expectedCount +\n      \", got \" + results.items.length + \"): \" + oldText.slice(0, 80)\n    );\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Updates the stack-trace text in the document (a Java exception dump\n# embedded as plain text in a single run) to reflect refactored line\n# numbers / stack frames, per a M2DocEvaluator 3.1.1 -> 3.2.0 code move:\n#  - several M2DocEvaluator.java / M2DocUtils.java / AbstractTemplatesTestSuite.java\n#    line numbers shift,\n#  - a generated accessor class name changes,\n#  - the tail of the trace (Eclipse JDT JUnit launcher frames) is replaced\n#    by a Maven Surefire / Tycho / Equinox launcher frame sequence.\n\n$d = $word.ActiveDocument\n\n$TAB = [string][char]9\n$NL  = [string][char]10\n\n# Each entry: old text to find, new text to replace with, expected match count.\n# Multi-line trace frames are tab-prefixed; built from line arrays joined with\n# ($NL + $TAB) so the real tab/newline bytes land in the string. Using single-\n# quoted PowerShell literals for every line means no backtick escapes are needed\n# and the literal dollar sign in 'ProviderFactory$ProviderProxy' is never treated\n# as variable interpolation.\n$replacements = @()\n\n$replacements += [PSCustomObject]@{ Old = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1207)'; New = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1296)'; Count = 1 }\n$replacements += [PSCustomObject]@{ Old = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)'; New = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)'; Count = 3 }\n$replacements += [PSCustomObject]@{ Old = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1467)'; New = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1556)'; Count = 1 }\n$replacements += [PSCustomObject]@{ Old = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:297)'; New = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:301)'; Count = 1 }\n$replacements += [PSCustomObject]@{ Old = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:282)'; New = 'org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:286)'; Count = 1 }\n$replacements += [PSCustomObject]@{ Old = 'org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:845)'; New = 'org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)'; Count = 1 }\n$replacements += [PSCustomObject]@{ Old = 'org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514)'; New = 'org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)'; Count = 1 }\n$replacements += [PSCustomObject]@{ Old = 'org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421)'; New = 'org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)'; Count = 1 }\n$replacements += [PSCustomObject]@{ Old = 'sun.reflect.GeneratedMethodAccessor73.invoke(Unknown Source)'; New = 'sun.reflect.GeneratedMethodAccessor5.invoke(Unknown Source)'; Count = 1 }\n$oldLines = @(\n  'at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)',\n  'at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)',\n  'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)',\n  'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)',\n  'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)',\n  'at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)'\n)\n$newLines = @(\n  'at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)',\n  'at org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)',\n  'at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)',\n  'at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)',\n  'at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)',\n  'at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)',\n  'at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)',\n  'at java.lang.reflect.Method.invoke(Method.java:498)',\n  'at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)',\n  'at org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:161)',\n  'at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)',\n  'at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)',\n  'at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)',\n  'at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)',\n  'at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)',\n  'at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)',\n  'at java.lang.reflect.Method.invoke(Method.java:498)',\n  'at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)',\n  'at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)',\n  'at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)',\n  'at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)',\n  'at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)',\n  'at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)',\n  'at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)',\n  'at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)',\n  'at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)',\n  'at java.lang.reflect.Method.invoke(Method.java:498)',\n  'at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)',\n  'at org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)',\n  'at org.eclipse.equinox.launcher.Main.run(Main.java:1447)',\n  'at org.eclipse.equinox.launcher.Main.main(Main.java:1420)'\n)\n$oldBlock = $TAB + ($oldLines -join ($NL + $TAB))\n$newBlock = $TAB + ($newLines -join ($NL + $TAB))\n$replacements += [PSCustomObject]@{ Old = $oldBlock; New = $newBlock; Count = 1 }\n\nforeach ($r in $replacements) {\n    # Sanity-check the expected number of occurrences before touching anything,\n    # using a plain-text regex scan over the whole story (robust to the old\n    # text spanning tabs/newlines inside a single run).\n    $beforeText = $d.Content.Text\n    $preCount = ([regex]::Matches($beforeText, [regex]::Escape($r.Old))).Count\n    if ($preCount -ne $r.Count) {\n        throw \"Unexpected pre-replace match count (expected $($r.Count), got $preCount) for: $($r.Old.Substring(0, [Math]::Min(80, $r.Old.Length)))\"\n    }\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 0            # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2) | Out-Null\n\n    # Confirm the substitution happened the expected number of times.\n    $afterText = $d.Content.Text\n    $postCount = ([regex]::Matches($afterText, [regex]::Escape($r.New))).Count\n    if ($postCount -lt $r.Count) {\n        throw \"Replacement did not apply as expected (expected >= $($r.Count), got $postCount) for: $($r.New.Substring(0, [Math]::Min(80, $r.New.Length)))\"\n    }\n}\n"}
